# New crime data collected — weekly refresh of the 115th Precinct CompStat
# report: bump the "Volume/Number" and reporting-week header text, then
# push the refreshed weekly/28-day/YTD/2-year/13-year/30-year figures into
# the data table (rows 15-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text cells) - edit only the runs that actually
# changed, leaving the rest of each cell's formatting untouched.
# ---------------------------------------------------------------------

# A8: "Volume 30   Number  12" -> "...13"
$ws.Range("A8").Characters(21, 2).Text = "13"

# C9: "Report Covering the Week  3/20/2023  Through  3/26/2023"
#  -> "Report Covering the Week  3/27/2023  Through  4/2/2023"
$ws.Range("C9").Characters(27, 9).Text = "3/27/2023"
$ws.Range("C9").Characters(47, 9).Text = "4/2/2023"

# ---------------------------------------------------------------------
# Row 15 (Murder)
# ---------------------------------------------------------------------
$ws.Range("C15").Value = "'0"
$ws.Range("E15").Value = -100
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -14.285714285714
$ws.Range("L15").Value = -33.333333333333
$ws.Range("M15").Value = 20

# ---------------------------------------------------------------------
# Row 16 (Rape)
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 61
$ws.Range("J16").Value = 65
$ws.Range("K16").Value = -6.153846153846
$ws.Range("L16").Value = 15.094339622641
$ws.Range("M16").Value = -19.736842105263
$ws.Range("N16").Value = -81.402439024390

# ---------------------------------------------------------------------
# Row 17 (Robbery)
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 37
$ws.Range("H17").Value = 5.714285714285
$ws.Range("I17").Value = 119
$ws.Range("J17").Value = 101
$ws.Range("K17").Value = 17.821782178217
$ws.Range("L17").Value = 67.605633802816
$ws.Range("M17").Value = 67.605633802816
$ws.Range("N17").Value = 52.564102564102

# ---------------------------------------------------------------------
# Row 18 (Fel. Assault)
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 22.222222222222
$ws.Range("I18").Value = 28
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = -22.222222222222
$ws.Range("L18").Value = -20
$ws.Range("M18").Value = -64.102564102564
$ws.Range("N18").Value = -95.483870967741

# ---------------------------------------------------------------------
# Row 19 (Burglary)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 7.692307692307
$ws.Range("F19").Value = 69
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = 32.692307692307
$ws.Range("I19").Value = 198
$ws.Range("J19").Value = 288
$ws.Range("K19").Value = -31.25
$ws.Range("L19").Value = 51.145038167938
$ws.Range("M19").Value = 88.571428571428
$ws.Range("N19").Value = -34.868421052631

# ---------------------------------------------------------------------
# Row 20 (Gr. Larceny)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 9.090909090909
$ws.Range("I20").Value = 78
$ws.Range("J20").Value = 61
$ws.Range("K20").Value = 27.868852459016
$ws.Range("L20").Value = 122.857142857143
$ws.Range("M20").Value = 32.203389830508
$ws.Range("N20").Value = -85.869565217391

# ---------------------------------------------------------------------
# Row 21 (G.L.A.)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = 8.108108108108
$ws.Range("F21").Value = 161
$ws.Range("G21").Value = 136
$ws.Range("H21").Value = 18.382352941176
$ws.Range("I21").Value = 490
$ws.Range("J21").Value = 560
$ws.Range("K21").Value = -12.5
$ws.Range("L21").Value = 46.268656716417
$ws.Range("M21").Value = 24.050632911392
$ws.Range("N21").Value = -74.169741697417

# ---------------------------------------------------------------------
# Row 22 (TOTAL)
# ---------------------------------------------------------------------
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("F22").Value = 13
$ws.Range("H22").Value = 550
$ws.Range("I22").Value = 28
$ws.Range("K22").Value = 180
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = 211.111111111111

# ---------------------------------------------------------------------
# Row 24 (Housing)
# ---------------------------------------------------------------------
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 84.210526315789
$ws.Range("F24").Value = 151
$ws.Range("G24").Value = 113
$ws.Range("H24").Value = 33.628318584070
$ws.Range("I24").Value = 554
$ws.Range("J24").Value = 413
$ws.Range("K24").Value = 34.140435835351
$ws.Range("L24").Value = 84.666666666666
$ws.Range("M24").Value = 134.745762711864

# ---------------------------------------------------------------------
# Row 25 (Petit Larceny)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 66.666666666666
$ws.Range("F25").Value = 73
$ws.Range("G25").Value = 52
$ws.Range("H25").Value = 40.384615384615
$ws.Range("I25").Value = 217
$ws.Range("J25").Value = 221
$ws.Range("K25").Value = -1.809954751131
$ws.Range("L25").Value = 33.950617283950
$ws.Range("M25").Value = -6.060606060606

# ---------------------------------------------------------------------
# Row 26 (Misd. Assault)
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 33.333333333333
$ws.Range("I26").Value = 13
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = 18.181818181818

# ---------------------------------------------------------------------
# Row 27 (UCR Rape*)
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 11
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 120
$ws.Range("I27").Value = 27
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = 22.727272727272
$ws.Range("L27").Value = -12.903225806451

# ---------------------------------------------------------------------
# Row 28 (Other Sex Crimes) - only the 30-year % chg moved
# ---------------------------------------------------------------------
$ws.Range("N28").Value = -92.857142857142

# ---------------------------------------------------------------------
# Row 29 (Shooting Vic.) - only the 30-year % chg moved
# ---------------------------------------------------------------------
$ws.Range("N29").Value = -92.307692307692

# ---------------------------------------------------------------------
# Row 30 (Shooting Inc.)
# ---------------------------------------------------------------------
$ws.Range("I30").Value = 3
$ws.Range("K30").Value = 200
$ws.Range("L30").Value = 200
$ws.Range("L30").NumberFormat = '#,##0.0;"-"#,##0.0'
